$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted right before the current
# row 41, pushing the existing rows 41-48 down to 42-49.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new record's data.
$ws.Range("A41").Value = 1
$ws.Range("B41").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C41").Value = "Arica y Parinacota"
$ws.Range("D41").Value = 45154
$ws.Range("E41").Value = 15
$ws.Range("F41").Value = 100112003
$ws.Range("G41").Value = "Ajo"
$ws.Range("H41").Value = "Chino"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 200
$ws.Range("K41").Value = 17000
$ws.Range("L41").Value = 18000
$ws.Range("M41").Value = 17500
$ws.Range("N41").Value = "$/caja 10 kilos"
$ws.Range("O41").Value = "China"
$ws.Range("P41").Value = 1750
$ws.Range("Q41").Value = 10
$ws.Range("R41").Value = "Hortaliza"
